$wb = $excel.ActiveWorkbook

# --- Update selection on the function_parameters sheet (was tabSelected) ---
$wsFuncParams = $wb.Worksheets.Item("function_parameters")
$wsFuncParams.Range("G6").Select()

# --- Add the two new sheets at the end of the workbook ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)

$ws1 = $wb.Worksheets.Add($null, $lastSheet)
$ws1.Name = "TreeStructre"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Trees"

# --- TreeStructre sheet data ---
$ws1.Cells.Item(1,1).Value = "id"
$ws1.Cells.Item(1,2).Value = "tree_id"
$ws1.Cells.Item(1,3).Value = "node_id"
$ws1.Cells.Item(1,4).Value = "node_name"
$ws1.Cells.Item(1,5).Value = "node_data"
$ws1.Cells.Item(1,6).Value = "parent"
$ws1.Cells.Item(2,1).Value = 1
$ws1.Cells.Item(2,2).Value = 1
$ws1.Cells.Item(2,3).Value = 1
$ws1.Cells.Item(2,4).Value = "here"
$ws1.Cells.Item(2,5).Value = "once"
$ws1.Cells.Item(2,6).Value = 0
$ws1.Cells.Item(3,1).Value = 2
$ws1.Cells.Item(3,2).Value = 1
$ws1.Cells.Item(3,3).Value = 2
$ws1.Cells.Item(3,4).Value = "is"
$ws1.Cells.Item(3,5).Value = "upon a time"
$ws1.Cells.Item(3,6).Value = 1
$ws1.Cells.Item(4,1).Value = 3
$ws1.Cells.Item(4,2).Value = 1
$ws1.Cells.Item(4,3).Value = 3
$ws1.Cells.Item(4,4).Value = "some"
$ws1.Cells.Item(4,5).Value = "there was"
$ws1.Cells.Item(4,6).Value = 1
$ws1.Cells.Item(5,1).Value = 4
$ws1.Cells.Item(5,2).Value = 1
$ws1.Cells.Item(5,3).Value = 4
$ws1.Cells.Item(5,4).Value = "data"
$ws1.Cells.Item(5,5).Value = "a little"
$ws1.Cells.Item(5,6).Value = 1
$ws1.Cells.Item(6,1).Value = 5
$ws1.Cells.Item(6,2).Value = 1
$ws1.Cells.Item(6,3).Value = 5
$ws1.Cells.Item(6,4).Value = "I"
$ws1.Cells.Item(6,5).Value = "girl name"
$ws1.Cells.Item(6,6).Value = 3
$ws1.Cells.Item(7,1).Value = 6
$ws1.Cells.Item(7,2).Value = 1
$ws1.Cells.Item(7,3).Value = 6
$ws1.Cells.Item(7,4).Value = "really"
$ws1.Cells.Item(7,5).Value = "snowwhite"
$ws1.Cells.Item(7,6).Value = 4
$ws1.Cells.Item(8,1).Value = 7
$ws1.Cells.Item(8,2).Value = 1
$ws1.Cells.Item(8,3).Value = 7
$ws1.Cells.Item(8,4).Value = "want"
$ws1.Cells.Item(8,5).Value = "dgdfg"
$ws1.Cells.Item(8,6).Value = 4
$ws1.Cells.Item(9,1).Value = 8
$ws1.Cells.Item(9,2).Value = 1
$ws1.Cells.Item(9,3).Value = 8
$ws1.Cells.Item(9,4).Value = "you"
$ws1.Cells.Item(9,5).Value = "dfgfd"
$ws1.Cells.Item(9,6).Value = 4
$ws1.Cells.Item(10,1).Value = 9
$ws1.Cells.Item(10,2).Value = 1
$ws1.Cells.Item(10,3).Value = 9
$ws1.Cells.Item(10,4).Value = "to"
$ws1.Cells.Item(10,5).Value = "ddd"
$ws1.Cells.Item(10,6).Value = 5
$ws1.Cells.Item(11,1).Value = 10
$ws1.Cells.Item(11,2).Value = 1
$ws1.Cells.Item(11,3).Value = 10
$ws1.Cells.Item(11,4).Value = "save"
$ws1.Cells.Item(11,5).Value = "ds"
$ws1.Cells.Item(11,6).Value = 7
$ws1.Cells.Item(12,1).Value = 11
$ws1.Cells.Item(12,2).Value = 2
$ws1.Cells.Item(12,3).Value = 1
$ws1.Cells.Item(12,4).Value = "some"
$ws1.Cells.Item(12,5).Value = "aasssad"
$ws1.Cells.Item(12,6).Value = 0
$ws1.Cells.Item(13,1).Value = 12
$ws1.Cells.Item(13,2).Value = 2
$ws1.Cells.Item(13,3).Value = 2
$ws1.Cells.Item(13,4).Value = "simple"
$ws1.Cells.Item(13,5).Value = "asdsad"
$ws1.Cells.Item(13,6).Value = 1
$ws1.Cells.Item(14,1).Value = 13
$ws1.Cells.Item(14,2).Value = 2
$ws1.Cells.Item(14,3).Value = 3
$ws1.Cells.Item(14,4).Value = "tree"
$ws1.Cells.Item(14,5).Value = "ddd"
$ws1.Cells.Item(14,6).Value = 1
$ws1.Cells.Item(15,1).Value = 14
$ws1.Cells.Item(15,2).Value = 3
$ws1.Cells.Item(15,3).Value = 1
$ws1.Cells.Item(15,4).Value = "another"
$ws1.Cells.Item(15,5).Value = "daeee"
$ws1.Cells.Item(15,6).Value = 0
$ws1.Cells.Item(16,1).Value = 15
$ws1.Cells.Item(16,2).Value = 3
$ws1.Cells.Item(16,3).Value = 2
$ws1.Cells.Item(16,4).Value = "simple"
$ws1.Cells.Item(16,5).Value = "rrrrr"
$ws1.Cells.Item(16,6).Value = 1
$ws1.Cells.Item(17,1).Value = 16
$ws1.Cells.Item(17,2).Value = 3
$ws1.Cells.Item(17,3).Value = 3
$ws1.Cells.Item(17,4).Value = "tree"
$ws1.Cells.Item(17,5).Value = "ttt"
$ws1.Cells.Item(17,6).Value = 2
$ws1.Cells.Item(18,1).Value = 17
$ws1.Cells.Item(18,2).Value = 4
$ws1.Cells.Item(18,3).Value = 1
$ws1.Cells.Item(18,4).Value = "oneandonly"
$ws1.Cells.Item(18,5).Value = "ffffffs"
$ws1.Cells.Item(18,6).Value = 0
# --- Trees sheet data ---
$ws2.Cells.Item(1,1).Value = "id"
$ws2.Cells.Item(1,2).Value = "name"
$ws2.Cells.Item(1,3).Value = "function"
$ws2.Cells.Item(1,4).Value = "nodes"
$ws2.Cells.Item(2,1).Value = 1
$ws2.Cells.Item(2,2).Value = "big tree"
$ws2.Cells.Item(2,3).Value = 1
$ws2.Cells.Item(2,4).Value = 0
$ws2.Cells.Item(3,1).Value = 2
$ws2.Cells.Item(3,2).Value = "small tree"
$ws2.Cells.Item(3,3).Value = 2
$ws2.Cells.Item(3,4).Value = 0
$ws2.Cells.Item(4,1).Value = 3
$ws2.Cells.Item(4,2).Value = "small tree2"
$ws2.Cells.Item(4,3).Value = 3
$ws2.Cells.Item(4,4).Value = 0
$ws2.Cells.Item(5,1).Value = 4
$ws2.Cells.Item(5,2).Value = "null_tree"
$ws2.Cells.Item(5,3).Value = 4
$ws2.Cells.Item(5,4).Value = 0
# --- Set the selection on the Trees sheet ---
$ws2.Range("C5").Select()

# --- Activate TreeStructre (becomes the workbook's active/selected tab) and set its selection ---
$ws1.Activate()
$ws1.Range("E1").Select()
